$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so Excel stores them as text (matching original inlineStr string cells)
# instead of auto-converting to numeric values.
$textForceCells = @("D5","D6","D7","D10","D12","D13","D14","D17","D19","D22","D23","D28","D29","D30","D31","D32","D33","D34","D35","D36","D46","D47","D48","D49","D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "43.039.32"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3
$ws.Range("D3").Value = "2.307.42"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "300.32"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6
$ws.Range("D6").Value = "97.93"
$ws.Range("E6").Value = "  -0.59%  "

# Row 7
$ws.Range("D7").Value = "0.513"
$ws.Range("E7").Value = "  -2.16%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  -2.84%  "

# Row 10
$ws.Range("D10").Value = "36.00"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("E11").Value = "  +0.09%  "

# Row 12
$ws.Range("D12").Value = "18.17"
$ws.Range("E12").Value = "  +1.21%  "

# Row 13
$ws.Range("D13").Value = "0.118"
$ws.Range("E13").Value = "  +1.82%  "

# Row 14
$ws.Range("D14").Value = "6.81"
$ws.Range("E14").Value = "  -1.22%  "

# Row 15
$ws.Range("D15").Value = "2.667.04"
$ws.Range("E15").Value = "  +0.16%  "

# Row 16
$ws.Range("D16").Value = "2.297.47"
$ws.Range("E16").Value = "  +1.05%  "

# Row 17
$ws.Range("D17").Value = "0.781"
$ws.Range("E17").Value = "  -0.97%  "

# Row 18
$ws.Range("D18").Value = "42.984.88"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19
$ws.Range("D19").Value = "12.71"
$ws.Range("E19").Value = "  -5.01%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -0.37%  "

# Row 21
$ws.Range("E21").Value = "  -1.87%  "

# Row 22
$ws.Range("D22").Value = "68.12"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").Value = "240.68"
$ws.Range("E23").Value = "  +0.50%  "

# Row 24
$ws.Range("E24").Value = "  -0.66%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("E27").Value = "  +0.21%  "

# Row 28
$ws.Range("D28").Value = "25.52"
$ws.Range("E28").Value = "  +2.95%  "

# Row 29
$ws.Range("D29").Value = "165.77"
$ws.Range("E29").Value = "  -1.41%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.08"
$ws.Range("E30").Value = "  -0.57%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -0.36%  "

# Row 32
$ws.Range("D32").Value = "33.20"
$ws.Range("E32").Value = "  -0.51%  "

# Row 33
$ws.Range("D33").Value = "4.96"
$ws.Range("E33").Value = "  +2.57%  "

# Row 34
$ws.Range("D34").Value = "1.00"

# Row 35
$ws.Range("D35").Value = "5.04"
$ws.Range("E35").Value = "  -3.61%  "

# Row 36
$ws.Range("D36").Value = "17.05"
$ws.Range("E36").Value = "  -6.39%  "

# Row 37
$ws.Range("E37").Value = "  -1.08%  "

# Row 38
$ws.Range("E38").Value = "  -0.44%  "

# Row 39
$ws.Range("E39").Value = "  -0.71%  "

# Row 40
$ws.Range("E40").Value = "  -1.48%  "

# Row 41
$ws.Range("E41").Value = "  -0.28%  "

# Row 42
$ws.Range("E42").Value = "  -1.60%  "

# Row 43
$ws.Range("D43").Value = "2.011.99"
$ws.Range("E43").Value = "  +0.70%  "

# Row 44
$ws.Range("E44").Value = "  -2.15%  "

# Row 45
$ws.Range("E45").Value = "  +1.79%  "

# Row 46
$ws.Range("D46").Value = "10.13"
$ws.Range("E46").Value = "  +0.30%  "

# Row 47
$ws.Range("D47").Value = "17.45"
$ws.Range("E47").Value = "  -0.58%  "

# Row 48
$ws.Range("D48").Value = "2.80"
$ws.Range("E48").Value = "  -0.87%  "

# Row 49
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").Value = "  -3.29%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "53.77"
$ws.Range("E50").Value = "  -1.54%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.533.76"
$ws.Range("E51").Value = "  +0.16%  "
